$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1266.3334
$ws.Range("I2").Value = 525
$ws.Range("J2").Value = 2749
$ws.Range("K2").Value = 525
$ws.Range("L2").Value = 2749
$ws.Range("M2").Value = -412
$ws.Range("N2").Value = -2975

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1793.9231
$ws.Range("I4").Value = 1564.8
$ws.Range("K4").Value = 1564.8
$ws.Range("M4").Value = -1450.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 577.1667
$ws.Range("I9").Value = 520.2857
$ws.Range("J9").Value = 656.8
$ws.Range("K9").Value = 520.2857
$ws.Range("L9").Value = 656.8
$ws.Range("M9").Value = -351.2857
$ws.Range("N9").Value = -994.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 168333.17
$ws.Range("J17").Value = 168333.17
$ws.Range("L17").Value = 504999.51
$ws.Range("N17").Value = -505335.51

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 30953808
$ws.Range("I19").Value = 55556684
$ws.Range("J19").Value = 12501653
$ws.Range("K19").Value = 55556684
$ws.Range("L19").Value = 12501653
$ws.Range("M19").Value = -55556509
$ws.Range("N19").Value = -12502003

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 3069.5386
$ws.Range("I41").Value = 3914.25
$ws.Range("J41").Value = 1718
$ws.Range("K41").Value = 3914.25
$ws.Range("L41").Value = 1718
$ws.Range("M41").Value = -3474.25
$ws.Range("N41").Value = -2598

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 351.93332
$ws.Range("I55").Value = 205.22223
$ws.Range("K55").Value = 205.22223
$ws.Range("M55").Value = 8.777770000000004

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 3050153
$ws.Range("I70").Value = 12195121
$ws.Range("J70").Value = 1830.3334
$ws.Range("K70").Value = 36585363
$ws.Range("L70").Value = 5491.0002
$ws.Range("M70").Value = -36585093
$ws.Range("N70").Value = -6031.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 3050153
$ws.Range("I73").Value = 12195121
$ws.Range("J73").Value = 1830.3334
$ws.Range("K73").Value = 36585363
$ws.Range("L73").Value = 5491.0002
$ws.Range("M73").Value = -36584427
$ws.Range("N73").Value = -7363.0002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 7014.2856
$ws.Range("I74").Value = 6516.6665
$ws.Range("K74").Value = 6516.6665
$ws.Range("M74").Value = -5580.6665

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 7014.2856
$ws.Range("I77").Value = 6516.6665
$ws.Range("K77").Value = 32583.3325
$ws.Range("M77").Value = -27903.3325

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 5102.625
$ws.Range("J112").Value = 5212.8
$ws.Range("L112").Value = 15638.4
$ws.Range("N112").Value = -17854.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2190.795
$ws.Range("I132").Value = 2209.7778
$ws.Range("K132").Value = 6629.3334
$ws.Range("M132").Value = -4099.3334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 5354.4204
$ws.Range("I138").Value = 3877.85
$ws.Range("J138").Value = 5957.102
$ws.Range("K138").Value = 11633.55
$ws.Range("L138").Value = 17871.306
$ws.Range("M138").Value = -6493.549999999999
$ws.Range("N138").Value = -28151.306

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2649.5557
$ws.Range("I141").Value = 2649.5557
$ws.Range("K141").Value = 7948.6671
$ws.Range("M141").Value = -2768.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1762
$ws.Range("I2").Value = 1750.5714
$ws.Range("J2").Value = 1778
$ws.Range("K2").Value = 1750.5714
$ws.Range("L2").Value = 1778
$ws.Range("M2").Value = -1637.5714
$ws.Range("N2").Value = -2004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8345.849
$ws.Range("I32").Value = 7613.951
$ws.Range("K32").Value = 7613.951
$ws.Range("M32").Value = -7326.951

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3604.5557
$ws.Range("I45").Value = 962.5
$ws.Range("K45").Value = 962.5
$ws.Range("M45").Value = -585.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 18131.455
$ws.Range("I46").Value = 8519
$ws.Range("K46").Value = 8519
$ws.Range("M46").Value = -8200

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6065854
$ws.Range("I61").Value = 8338591
$ws.Range("K61").Value = 8338591
$ws.Range("M61").Value = -8338379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3193.25
$ws.Range("I63").Value = 3669.2
$ws.Range("J63").Value = 2400
$ws.Range("K63").Value = 3669.2
$ws.Range("L63").Value = 2400
$ws.Range("M63").Value = -2983.2
$ws.Range("N63").Value = -3772

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3193.25
$ws.Range("I66").Value = 3669.2
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 18346
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = -14914
$ws.Range("N66").Value = -18864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H104").Value = 25000
$ws.Range("J104").Value = 25000
$ws.Range("L104").Value = 25000
$ws.Range("N104").Value = -31988

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5073.857
$ws.Range("I110").Value = 5356.3687
$ws.Range("J110").Value = 4477.4443
$ws.Range("K110").Value = 5356.3687
$ws.Range("L110").Value = 4477.4443
$ws.Range("M110").Value = -3311.3687
$ws.Range("N110").Value = -8567.4443

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1762
$ws.Range("I116").Value = 1750.5714
$ws.Range("J116").Value = 1778
$ws.Range("K116").Value = 1750.5714
$ws.Range("L116").Value = 1778
$ws.Range("M116").Value = 543.4286
$ws.Range("N116").Value = -6366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4783.5654
$ws.Range("I122").Value = 3613.2354
$ws.Range("K122").Value = 10839.7062
$ws.Range("M122").Value = -8389.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3965.35
$ws.Range("I132").Value = 3965.35
$ws.Range("K132").Value = 11896.05
$ws.Range("M132").Value = -9366.05

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6065854
$ws.Range("I136").Value = 8338591
$ws.Range("K136").Value = 25015773
$ws.Range("M136").Value = -25013223

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1762
$ws.Range("I3").Value = 1750.5714
$ws.Range("J3").Value = 1778
$ws.Range("K3").Value = 1750.5714
$ws.Range("L3").Value = 1778
$ws.Range("M3").Value = -1636.5714
$ws.Range("N3").Value = -2006

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1409.7858
$ws.Range("I80").Value = 1322.5
$ws.Range("J80").Value = 1458.2778
$ws.Range("K80").Value = 1322.5
$ws.Range("L80").Value = 1458.2778
$ws.Range("M80").Value = -324.5
$ws.Range("N80").Value = -3454.2778

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 1409.7858
$ws.Range("I83").Value = 1322.5
$ws.Range("J83").Value = 1458.2778
$ws.Range("K83").Value = 6612.5
$ws.Range("L83").Value = 7291.389
$ws.Range("M83").Value = -1620.5
$ws.Range("N83").Value = -17275.389

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1914.5161
$ws.Range("I94").Value = 2148.8635
$ws.Range("J94").Value = 1341.6666
$ws.Range("K94").Value = 2148.8635
$ws.Range("L94").Value = 1341.6666
$ws.Range("M94").Value = -1697.8635
$ws.Range("N94").Value = -2243.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 5144.75
$ws.Range("I107").Value = 5220.2
$ws.Range("K107").Value = 5220.2
$ws.Range("M107").Value = -3300.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 72600
$ws.Range("J132").Value = 72600
$ws.Range("L132").Value = 72600
$ws.Range("N132").Value = -82720

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2523.4722
$ws.Range("I134").Value = 2065.1155
$ws.Range("J134").Value = 3715.2
$ws.Range("K134").Value = 6195.3465
$ws.Range("L134").Value = 11145.6
$ws.Range("M134").Value = -3660.3465
$ws.Range("N134").Value = -16215.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 120000
$ws.Range("J140").Value = 120000
$ws.Range("L140").Value = 120000
$ws.Range("N140").Value = -130360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 22730212
$ws.Range("I31").Value = 25643604
$ws.Range("J31").Value = 5742.4
$ws.Range("K31").Value = 25643604
$ws.Range("L31").Value = 5742.4
$ws.Range("M31").Value = -25643309
$ws.Range("N31").Value = -6332.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 22730212
$ws.Range("I34").Value = 25643604
$ws.Range("J34").Value = 5742.4
$ws.Range("K34").Value = 25643604
$ws.Range("L34").Value = 5742.4
$ws.Range("M34").Value = -25643402
$ws.Range("N34").Value = -6146.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3892.4
$ws.Range("I86").Value = 3788
$ws.Range("J86").Value = 3996.8
$ws.Range("K86").Value = 3788
$ws.Range("L86").Value = 3996.8
$ws.Range("M86").Value = -2665
$ws.Range("N86").Value = -6242.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3892.4
$ws.Range("I89").Value = 3788
$ws.Range("J89").Value = 3996.8
$ws.Range("K89").Value = 18940
$ws.Range("L89").Value = 19984
$ws.Range("M89").Value = -13324
$ws.Range("N89").Value = -31216

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H96").Value = 24331
$ws.Range("J96").Value = 24331
$ws.Range("L96").Value = 24331
$ws.Range("N96").Value = -29823

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 54999
$ws.Range("J106").Value = 54999
$ws.Range("L106").Value = 54999
$ws.Range("N106").Value = -57523

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 44999.5
$ws.Range("J111").Value = 44999.5
$ws.Range("L111").Value = 44999.5
$ws.Range("N111").Value = -53179.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 544
$ws.Range("I132").Value = 285.2857
$ws.Range("J132").Value = 1449.5
$ws.Range("K132").Value = 855.8571000000001
$ws.Range("L132").Value = 4348.5
$ws.Range("M132").Value = 1674.1429
$ws.Range("N132").Value = -9408.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2022.3429
$ws.Range("I134").Value = 2058.138
$ws.Range("J134").Value = 1849.3334
$ws.Range("K134").Value = 6174.414
$ws.Range("L134").Value = 5548.0002
$ws.Range("M134").Value = -3639.414
$ws.Range("N134").Value = -10618.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 385828.5
$ws.Range("J141").Value = 502712.72
$ws.Range("L141").Value = 502712.72
$ws.Range("N141").Value = -513072.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 988.0476
$ws.Range("I5").Value = 738.3571
$ws.Range("J5").Value = 1487.4286
$ws.Range("K5").Value = 2215.0713
$ws.Range("L5").Value = 4462.2858
$ws.Range("M5").Value = -2103.0713
$ws.Range("N5").Value = -4686.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 3288
$ws.Range("J29").Value = 3288
$ws.Range("L29").Value = 9864
$ws.Range("N29").Value = -10418

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 333335680
$ws.Range("I42").Value = 500001500
$ws.Range("K42").Value = 1500004500
$ws.Range("M42").Value = -1500003966

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 15163.833
$ws.Range("I51").Value = 7750.75
$ws.Range("J51").Value = 29990
$ws.Range("K51").Value = 23252.25
$ws.Range("L51").Value = 89970
$ws.Range("M51").Value = -22792.25
$ws.Range("N51").Value = -90890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 5237.125
$ws.Range("I60").Value = 333.33334
$ws.Range("J60").Value = 8179.4
$ws.Range("K60").Value = 1000.00002
$ws.Range("L60").Value = 24538.2
$ws.Range("M60").Value = -749.0000200000001
$ws.Range("N60").Value = -25040.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2387.25
$ws.Range("I81").Value = 1516.3334
$ws.Range("K81").Value = 4549.0002
$ws.Range("M81").Value = -3426.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 2387.25
$ws.Range("I84").Value = 1516.3334
$ws.Range("K84").Value = 13647.0006
$ws.Range("M84").Value = -8031.000599999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1752.1578
$ws.Range("I113").Value = 1988.6666
$ws.Range("K113").Value = 5965.9998
$ws.Range("M113").Value = -3795.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 5320.647
$ws.Range("I114").Value = 6323
$ws.Range("K114").Value = 18969
$ws.Range("M114").Value = -15715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 55914.832
$ws.Range("I122").Value = 110696.664
$ws.Range("J122").Value = 1133
$ws.Range("K122").Value = 996269.976
$ws.Range("L122").Value = 10197
$ws.Range("M122").Value = -993819.976
$ws.Range("N122").Value = -15097

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2965.0454
$ws.Range("I129").Value = 2256
$ws.Range("J129").Value = 3674.0908
$ws.Range("K129").Value = 6768
$ws.Range("L129").Value = 11022.2724
$ws.Range("M129").Value = -1768
$ws.Range("N129").Value = -21022.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2915.077
$ws.Range("J131").Value = 3131.0942
$ws.Range("L131").Value = 9393.2826
$ws.Range("N131").Value = -19473.2826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 11458.728
$ws.Range("I134").Value = 4505.75
$ws.Range("K134").Value = 13517.25
$ws.Range("M134").Value = -8447.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 988.0476
$ws.Range("I135").Value = 738.3571
$ws.Range("J135").Value = 1487.4286
$ws.Range("K135").Value = 6645.2139
$ws.Range("L135").Value = 13386.8574
$ws.Range("M135").Value = -4110.2139
$ws.Range("N135").Value = -18456.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3585.5789
$ws.Range("I140").Value = 1170.4166
$ws.Range("K140").Value = 3511.2498
$ws.Range("M140").Value = 1668.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 170.07143
$ws.Range("I2").Value = 119.333336
$ws.Range("K2").Value = 119.333336
$ws.Range("M2").Value = -6.333336000000003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 833.3333
$ws.Range("I12").Value = 750
$ws.Range("J12").Value = 1000
$ws.Range("K12").Value = 750
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = -610
$ws.Range("N12").Value = -1280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 3343332.8
$ws.Range("I14").Value = 14999
$ws.Range("J14").Value = 10000000
$ws.Range("K14").Value = 14999
$ws.Range("L14").Value = 10000000
$ws.Range("M14").Value = -14831
$ws.Range("N14").Value = -10000336

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 44444
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5526.875
$ws.Range("I80").Value = 3118.1667
$ws.Range("K80").Value = 3118.1667
$ws.Range("M80").Value = -2120.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 5526.875
$ws.Range("I83").Value = 3118.1667
$ws.Range("K83").Value = 15590.8335
$ws.Range("M83").Value = -10598.8335

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4148
$ws.Range("I102").Value = 3769.5557
$ws.Range("K102").Value = 3769.5557
$ws.Range("M102").Value = -2147.5557

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 281.4
$ws.Range("I107").Value = 203.33333
$ws.Range("J107").Value = 398.5
$ws.Range("K107").Value = 203.33333
$ws.Range("L107").Value = 398.5
$ws.Range("M107").Value = 1716.66667
$ws.Range("N107").Value = -4238.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 14180757
$ws.Range("I126").Value = 20381650
$ws.Range("J126").Value = 7290.7144
$ws.Range("K126").Value = 61144950
$ws.Range("L126").Value = 21872.1432
$ws.Range("M126").Value = -61142480
$ws.Range("N126").Value = -26812.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1928.9131
$ws.Range("I16").Value = 1928.9131
$ws.Range("K16").Value = 1928.9131
$ws.Range("M16").Value = -1758.9131

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 22004566
$ws.Range("I22").Value = 33004500
$ws.Range("J22").Value = 4700
$ws.Range("K22").Value = 33004500
$ws.Range("L22").Value = 4700
$ws.Range("M22").Value = -33004205
$ws.Range("N22").Value = -5290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 22004566
$ws.Range("I27").Value = 33004500
$ws.Range("J27").Value = 4700
$ws.Range("K27").Value = 33004500
$ws.Range("L27").Value = 4700
$ws.Range("M27").Value = -33004393
$ws.Range("N27").Value = -4914

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1119.4
$ws.Range("J46").Value = 1533
$ws.Range("L46").Value = 1533
$ws.Range("N46").Value = -1909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 865.64
$ws.Range("I55").Value = 561.6667
$ws.Range("J55").Value = 1146.2307
$ws.Range("K55").Value = 561.6667
$ws.Range("L55").Value = 1146.2307
$ws.Range("M55").Value = -388.6667
$ws.Range("N55").Value = -1492.2307

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2755.32
$ws.Range("I61").Value = 2738.3914
$ws.Range("J61").Value = 2950
$ws.Range("K61").Value = 2738.3914
$ws.Range("L61").Value = 2950
$ws.Range("M61").Value = -2536.3914
$ws.Range("N61").Value = -3354

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5002.8823
$ws.Range("I82").Value = 1272.6666
$ws.Range("J82").Value = 9199.375
$ws.Range("K82").Value = 1272.6666
$ws.Range("L82").Value = 9199.375
$ws.Range("M82").Value = -911.6666
$ws.Range("N82").Value = -9921.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 5002.8823
$ws.Range("I85").Value = 1272.6666
$ws.Range("J85").Value = 9199.375
$ws.Range("K85").Value = 1272.6666
$ws.Range("L85").Value = 9199.375
$ws.Range("M85").Value = -24.66660000000002
$ws.Range("N85").Value = -11695.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2059120.9
$ws.Range("I93").Value = 673
$ws.Range("J93").Value = 5558482.5
$ws.Range("K93").Value = 673
$ws.Range("L93").Value = 5558482.5
$ws.Range("M93").Value = 575
$ws.Range("N93").Value = -5560978.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2755.32
$ws.Range("I113").Value = 2738.3914
$ws.Range("J113").Value = 2950
$ws.Range("K113").Value = 2738.3914
$ws.Range("L113").Value = 2950
$ws.Range("M113").Value = -568.3914
$ws.Range("N113").Value = -7290

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3199.8154
$ws.Range("I122").Value = 3191.1091
$ws.Range("J122").Value = 3247.7
$ws.Range("K122").Value = 9573.3273
$ws.Range("L122").Value = 9743.099999999999
$ws.Range("M122").Value = -7123.327300000001
$ws.Range("N122").Value = -14643.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H124").Value = 69439.8
$ws.Range("J124").Value = 69439.8
$ws.Range("L124").Value = 69439.8
$ws.Range("N124").Value = -79259.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4117.84
$ws.Range("I132").Value = 3147.353
$ws.Range("J132").Value = 6180.125
$ws.Range("K132").Value = 9442.059000000001
$ws.Range("L132").Value = 18540.375
$ws.Range("M132").Value = -6912.059000000001
$ws.Range("N132").Value = -23600.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3690.5151
$ws.Range("I136").Value = 3561.724
$ws.Range("J136").Value = 4624.25
$ws.Range("K136").Value = 10685.172
$ws.Range("L136").Value = 13872.75
$ws.Range("M136").Value = -8135.172
$ws.Range("N136").Value = -18972.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1080.8
$ws.Range("J14").Value = 1400
$ws.Range("L14").Value = 1400
$ws.Range("N14").Value = -1736

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1100
$ws.Range("I17").Value = 1000
$ws.Range("K17").Value = 1000
$ws.Range("M17").Value = -828

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5183.9
$ws.Range("I81").Value = 2998.3333
$ws.Range("J81").Value = 8462.25
$ws.Range("K81").Value = 5996.6666
$ws.Range("L81").Value = 16924.5
$ws.Range("M81").Value = -4935.6666
$ws.Range("N81").Value = -19046.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 5183.9
$ws.Range("I84").Value = 2998.3333
$ws.Range("J84").Value = 8462.25
$ws.Range("K84").Value = 29983.333
$ws.Range("L84").Value = 84622.5
$ws.Range("M84").Value = -24679.333
$ws.Range("N84").Value = -95230.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 101966
$ws.Range("J88").Value = 101966
$ws.Range("L88").Value = 101966
$ws.Range("N88").Value = -102778

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H91").Value = 101966
$ws.Range("J91").Value = 101966
$ws.Range("L91").Value = 101966
$ws.Range("N91").Value = -104774

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 923.7
$ws.Range("I113").Value = 816.3
$ws.Range("K113").Value = 2448.9
$ws.Range("M113").Value = -278.8999999999996

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 49000
$ws.Range("J123").Value = 49000
$ws.Range("L123").Value = 49000
$ws.Range("N123").Value = -58800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2021.3704
$ws.Range("I132").Value = 1341.762
$ws.Range("J132").Value = 4400
$ws.Range("K132").Value = 4025.286
$ws.Range("L132").Value = 13200
$ws.Range("M132").Value = -1495.286
$ws.Range("N132").Value = -18260

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3502.653
$ws.Range("I136").Value = 3906.742
$ws.Range("K136").Value = 11720.226
$ws.Range("M136").Value = -9170.226
